$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "RQ004"
$ws.Range("B6").Value = "El sistema debe permitir upload de documentos"
$ws.Range("C6").Value = "Especificacion del cliente"
$ws.Range("E6").Value = "Caso de prueba upload recetas"
$ws.Range("D6").Value = "Disenio del sistema de almacenamiento"

$ws.Columns.Item(3).ColumnWidth = 30

[void]$ws.Range("B16").Select()
